$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple B-column value bumps for rows 2-10
$ws.Range("B2").Value = 80348
$ws.Range("B3").Value = 91808
$ws.Range("B4").Value = 91808
$ws.Range("B5").Value = 91808
$ws.Range("B6").Value = 80348
$ws.Range("B7").Value = 98902
$ws.Range("B8").Value = 80348
$ws.Range("B9").Value = 92179
$ws.Range("B10").Value = 80348

# Rows 11 and 12 effectively swap their content (A, D, E, F, G, H, Q, R),
# while column B gets independent new values.

# Target row 11 (after edit):
$ws.Range("A11").Value = 130894306
$ws.Range("B11").Value = 91771
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 5447
$ws.Range("F11").Value = "Vedticka"
$ws.Range("G11").Value = "Fuscoporia viticola"
$ws.Range("H11").Value = "(Schwein.) Murrill"
$ws.Range("Q11").Value = 799129
$ws.Range("R11").Value = 7351628

# Target row 12 (after edit):
$ws.Range("A12").Value = 130894267
$ws.Range("B12").Value = 91808
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 1202
$ws.Range("F12").Value = "Ullticka"
$ws.Range("G12").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H12").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 799039
$ws.Range("R12").Value = 7351523
